$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Power Detector row with the new part number and clear the
# (now unknown) footprint value.
$ws.Range("B5").Value = "AD8319"
$ws.Range("C5").ClearContents()

# Update the active selection to match the edited cell.
$ws.Range("C5").Select()
